$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-log entries (rows 6-10) appended below the existing data (rows 1-5).
# Values: Date (A), Start Time (B), End Time (C), Work Done (D)
$rows = @(
    @{ Row = 6;  Date = 45714; Start = 0.125;               End = 0.16666666666666666;  Work = "Worked on Overleaf Project Report" },
    @{ Row = 7;  Date = 45728; Start = 0.27777777777777779; End = 0.375;                 Work = "Worked on creating tasks and edited proposal for clearer reading" },
    @{ Row = 8;  Date = 45728; Start = 0.4375;               End = 0.47916666666666669;  Work = "Searched for Research papers related to project" },
    @{ Row = 9;  Date = 45729; Start = 0.16666666666666666; End = 0.20833333333333334;  Work = "Searched for Research papers related to project" },
    @{ Row = 10; Date = 45729; Start = 0.22916666666666666; End = 0.29166666666666669;  Work = "Worked on overleaf project report, planning section, updated project report with workflow" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the date/time number formats from row 4 (A4:C4) so the new cells
    # reuse the same cell styles instead of minting new ones.
    $ws.Range("A4:C4").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.Start
    $ws.Cells.Item($rowNum, 3).Value = $r.End
    $ws.Cells.Item($rowNum, 4).Value = $r.Work
}

# Leave the selection where the author left it when the file was saved.
$ws.Range("L14").Select() | Out-Null
